# pTHg_CombinedR.xlsx update
# Adds new rows of Combined Outflow Q/WQ data (16 -> 34 observations),
# refreshes the sheet selection/scroll position, and matches formatting
# of the new text-valued time cell (row 10, column B -> "14:20").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlRight = -4152

# --- Data rows (A:D), rows 3-36 -------------------------------------------------
$ws.Cells.Item(3,1).Value2 = 42389
$ws.Cells.Item(3,2).Value2 = 0.5625
$ws.Cells.Item(3,3).Value2 = 476
$ws.Cells.Item(3,4).Value2 = 33
$ws.Cells.Item(4,1).Value2 = 42392
$ws.Cells.Item(4,2).Value2 = 0.45833333333333331
$ws.Cells.Item(4,3).Value2 = 56
$ws.Cells.Item(4,4).Value2 = 12
$ws.Cells.Item(5,1).Value2 = 42393
$ws.Cells.Item(5,2).Value2 = 0.40183238440735369
$ws.Cells.Item(5,3).Value2 = 117
$ws.Cells.Item(5,4).Value2 = 8.4
$ws.Cells.Item(6,1).Value2 = 42402
$ws.Cells.Item(6,2).Value2 = 0.625
$ws.Cells.Item(6,3).Value2 = 498
$ws.Cells.Item(6,4).Value2 = 23
$ws.Cells.Item(7,1).Value2 = 42403
$ws.Cells.Item(7,2).Value2 = 0.58368055555555554
$ws.Cells.Item(7,3).Value2 = 430
$ws.Cells.Item(7,4).Value2 = 27
$ws.Cells.Item(8,1).Value2 = 42404
$ws.Cells.Item(8,2).Value2 = 0.52777777777777779
$ws.Cells.Item(8,3).Value2 = 256
$ws.Cells.Item(8,4).Value2 = 85
$ws.Cells.Item(9,1).Value2 = 42436
$ws.Cells.Item(9,2).Value2 = 0.56944444444444442
$ws.Cells.Item(9,3).Value2 = 541
$ws.Cells.Item(9,4).Value2 = 185
$ws.Cells.Item(10,1).Value2 = 42437
$ws.Cells.Item(10,2).Value = "14:20"
$ws.Cells.Item(10,2).HorizontalAlignment = $xlRight
$ws.Cells.Item(10,3).Value2 = 551
$ws.Cells.Item(10,4).Value2 = 188
$ws.Cells.Item(11,1).Value2 = 42439
$ws.Cells.Item(11,2).Value2 = 0.65933814661409929
$ws.Cells.Item(11,3).Value2 = 663
$ws.Cells.Item(11,4).Value2 = 31
$ws.Cells.Item(12,1).Value2 = 42441
$ws.Cells.Item(12,2).Value2 = 0.55815309248746492
$ws.Cells.Item(12,3).Value2 = 4378
$ws.Cells.Item(12,4).Value2 = 68
$ws.Cells.Item(13,1).Value2 = 42442
$ws.Cells.Item(13,2).Value2 = 0.51066760248567022
$ws.Cells.Item(13,3).Value2 = 3688
$ws.Cells.Item(13,4).Value2 = 55
$ws.Cells.Item(14,1).Value2 = 42443
$ws.Cells.Item(14,2).Value2 = 0.53119982396714727
$ws.Cells.Item(14,3).Value2 = 3190
$ws.Cells.Item(14,4).Value2 = 67
$ws.Cells.Item(15,1).Value2 = 42444
$ws.Cells.Item(15,2).Value2 = 0.51742016319439443
$ws.Cells.Item(15,3).Value2 = 2637
$ws.Cells.Item(15,4).Value2 = 42
$ws.Cells.Item(16,1).Value2 = 42446
$ws.Cells.Item(16,2).Value2 = 0.56646408071301035
$ws.Cells.Item(16,3).Value2 = 2840
$ws.Cells.Item(16,4).Value2 = 24
$ws.Cells.Item(17,1).Value2 = 42466
$ws.Cells.Item(17,2).Value2 = 0.54166666666666663
$ws.Cells.Item(17,3).Value2 = 473
$ws.Cells.Item(17,4).Value2 = 5.7
$ws.Cells.Item(18,1).Value2 = 42719
$ws.Cells.Item(18,2).Value2 = 0.69444444444444453
$ws.Cells.Item(18,3).Value2 = 79
$ws.Cells.Item(18,4).Value2 = 5.2
$ws.Cells.Item(19,1).Value2 = 42720
$ws.Cells.Item(19,2).Value2 = 0.72222222222222221
$ws.Cells.Item(19,3).Value2 = 233
$ws.Cells.Item(19,4).Value2 = 39
$ws.Cells.Item(20,1).Value2 = 42740
$ws.Cells.Item(20,2).Value2 = 0.61458333333333326
$ws.Cells.Item(20,3).Value2 = 205
$ws.Cells.Item(20,4).Value2 = 16
$ws.Cells.Item(21,1).Value2 = 42744
$ws.Cells.Item(21,2).Value2 = 0.72934110256754703
$ws.Cells.Item(21,3).Value2 = 12451
$ws.Cells.Item(21,4).Value2 = 272
$ws.Cells.Item(22,1).Value2 = 42745
$ws.Cells.Item(22,2).Value2 = 0.60417623205853732
$ws.Cells.Item(22,3).Value2 = 3801
$ws.Cells.Item(22,4).Value2 = 187
$ws.Cells.Item(23,1).Value2 = 42746
$ws.Cells.Item(23,2).Value2 = 0.66666817603807427
$ws.Cells.Item(23,3).Value2 = 13750
$ws.Cells.Item(23,4).Value2 = 161
$ws.Cells.Item(24,1).Value2 = 42749
$ws.Cells.Item(24,2).Value2 = 0.64930555555555558
$ws.Cells.Item(24,3).Value2 = 3891
$ws.Cells.Item(24,4).Value2 = 48
$ws.Cells.Item(25,1).Value2 = 42754
$ws.Cells.Item(25,2).Value2 = 0.52083333333333337
$ws.Cells.Item(25,3).Value2 = 6730
$ws.Cells.Item(25,4).Value2 = 63
$ws.Cells.Item(26,1).Value2 = 42755
$ws.Cells.Item(26,2).Value2 = 0.44791666666666669
$ws.Cells.Item(26,3).Value2 = 6730
$ws.Cells.Item(26,4).Value2 = 187
$ws.Cells.Item(27,1).Value2 = 42756
$ws.Cells.Item(27,2).Value2 = 0.3923611111111111
$ws.Cells.Item(27,3).Value2 = 10109
$ws.Cells.Item(27,4).Value2 = 180
$ws.Cells.Item(28,1).Value2 = 42758
$ws.Cells.Item(28,2).Value2 = 0.49652777777777779
$ws.Cells.Item(28,3).Value2 = 9896
$ws.Cells.Item(28,4).Value2 = 125
$ws.Cells.Item(29,1).Value2 = 42761
$ws.Cells.Item(29,2).Value2 = 0.47934027777777777
$ws.Cells.Item(29,3).Value2 = 3453
$ws.Cells.Item(29,4).Value2 = 19
$ws.Cells.Item(30,1).Value2 = 42767
$ws.Cells.Item(30,2).Value2 = 0.54185185925636103
$ws.Cells.Item(30,3).Value2 = 3053
$ws.Cells.Item(30,4).Value2 = 8.3000000000000007
$ws.Cells.Item(31,1).Value2 = 42774
$ws.Cells.Item(31,2).Value2 = 0.57639118479022766
$ws.Cells.Item(31,3).Value2 = 11060
$ws.Cells.Item(31,4).Value2 = 178
$ws.Cells.Item(32,1).Value2 = 42776
$ws.Cells.Item(32,2).Value2 = 0.61458485941221541
$ws.Cells.Item(32,3).Value2 = 11505
$ws.Cells.Item(32,4).Value2 = 115
$ws.Cells.Item(33,1).Value2 = 42787
$ws.Cells.Item(33,2).Value2 = 0.5590286977733987
$ws.Cells.Item(33,3).Value2 = 13460
$ws.Cells.Item(33,4).Value2 = 148
$ws.Cells.Item(34,1).Value2 = 42810
$ws.Cells.Item(34,2).Value2 = 0.64976728373393811
$ws.Cells.Item(34,3).Value2 = 2942
$ws.Cells.Item(34,4).Value2 = 4.9000000000000004
$ws.Cells.Item(35,1).Value2 = 42829
$ws.Cells.Item(35,2).Value2 = 0.65605302031566826
$ws.Cells.Item(35,3).Value2 = 259
$ws.Cells.Item(35,4).Value2 = 7.1
$ws.Cells.Item(36,1).Value2 = 42851
$ws.Cells.Item(36,2).Value2 = 0.60451388888888891
$ws.Cells.Item(36,3).Value2 = 223
$ws.Cells.Item(36,4).Value2 = 4.7

# --- Sheet view: scroll position + active selection ----------------------------
$ws.Range("E36").Select()

Write-Host "pTHg_CombinedR updated: 34 observations (rows 3-36)"
